$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 108 ("Rien ne nous concerne aujourd'hui !") had a stray "NA" value in
# its page-number column; clear it so it matches every other such row,
# which leaves that column blank.
$ws.Range("C108").Value = ""

# Append the new scraped result row for 2025-05-13, which inherits the
# "NA" page-number value that used to (incorrectly) sit on row 108.
# The date is entered with a leading apostrophe and the style is reset to
# "Normal" so Excel stores it as plain text ("2025-05-13"), not as a date
# serial number, matching how every other date cell in the sheet is saved.
$ws.Range("A109").Value = "'2025-05-13"
$ws.Range("A109").Style = "Normal"
$ws.Range("B109").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C109").Value = "NA"
$ws.Range("D109").Value = 1
